$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row (row 1) so the "_old"/"_new" suffixes used for the
# two compared AHB format versions become the concrete format-version names
# "_FV2404" and "_FV2410".
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the existing A1:U69 range into an actual Excel table ("Table1") that
# covers the whole sheet (header row + all 68 data rows), matching the
# newly-added xl/tables/table1.xml part and the sheet's <tableParts>.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U69"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (ySplit = 1, top-left cell of the scrolling area is
# A2), which is the <pane .../> addition seen in the sheetView.
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)
